$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.473279237747192
$ws.Range("B1").Value = 1.640435099601746
$ws.Range("C1").Value = 1.698710441589355
$ws.Range("D1").Value = 2.161072492599487
$ws.Range("E1").Value = 3.369428157806396
